$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 9.375560273003021
$ws.Cells.Item(2, 3).Value = 4.873334640101121
$ws.Cells.Item(2, 4).Value = 6.000359434813047
$ws.Cells.Item(2, 5).Value = 12.48587606959194
$ws.Cells.Item(2, 7).Value = 3.660882570932025
$ws.Cells.Item(2, 11).Value = 8.703586337095217
$ws.Cells.Item(2, 13).Value = 13.7879471518383
$ws.Cells.Item(2, 15).Value = 26.48908835048002

$ws.Cells.Item(3, 2).Value = 9.141207699389019
$ws.Cells.Item(3, 3).Value = 4.743231214474398
$ws.Cells.Item(3, 4).Value = 5.884439380121115
$ws.Cells.Item(3, 5).Value = 12.26020140850136
$ws.Cells.Item(3, 7).Value = 3.663165921118523
$ws.Cells.Item(3, 11).Value = 8.55653884931032
$ws.Cells.Item(3, 13).Value = 13.63813880662676
$ws.Cells.Item(3, 15).Value = 26.49545803006296

$ws.Cells.Item(4, 2).Value = 8.996669255938162
$ws.Cells.Item(4, 3).Value = 4.660468743089274
$ws.Cells.Item(4, 4).Value = 5.813895740404158
$ws.Cells.Item(4, 5).Value = 12.12370435956095
$ws.Cells.Item(4, 7).Value = 3.664641469367243
$ws.Cells.Item(4, 11).Value = 8.467074643385754
$ws.Cells.Item(4, 13).Value = 13.54902084635102
$ws.Cells.Item(4, 15).Value = 26.50491228906302

$ws.Cells.Item(5, 2).Value = 8.937703955803448
$ws.Cells.Item(5, 3).Value = 4.626044378126353
$ws.Cells.Item(5, 4).Value = 5.78535134227424
$ws.Cells.Item(5, 5).Value = 12.06868101070824
$ws.Cells.Item(5, 7).Value = 3.665261326821154
$ws.Cells.Item(5, 11).Value = 8.430875363724761
$ws.Cells.Item(5, 13).Value = 13.51346557096962
$ws.Cells.Item(5, 15).Value = 26.51015731744985

$ws.Cells.Item(6, 2).Value = 8.927911826222058
$ws.Cells.Item(6, 3).Value = 4.620286878125405
$ws.Cells.Item(6, 4).Value = 5.780625096237505
$ws.Cells.Item(6, 5).Value = 12.05958301541363
$ws.Cells.Item(6, 7).Value = 3.665365376448888
$ws.Cells.Item(6, 11).Value = 8.424881582985273
$ws.Cells.Item(6, 13).Value = 13.50760876254955
$ws.Cells.Item(6, 15).Value = 26.51111228605363

$ws.Cells.Item(7, 2).Value = 8.995874159718426
$ws.Cells.Item(7, 3).Value = 4.660007274501707
$ws.Cells.Item(7, 4).Value = 5.81350990176212
$ws.Cells.Item(7, 5).Value = 12.12295975922202
$ws.Cells.Item(7, 7).Value = 3.664649753753539
$ws.Cells.Item(7, 11).Value = 8.466585334679342
$ws.Cells.Item(7, 13).Value = 13.54853820172009
$ws.Cells.Item(7, 15).Value = 26.50497739041849

$ws.Cells.Item(8, 2).Value = 9.294946644677413
$ws.Cells.Item(8, 3).Value = 4.829087232544311
$ws.Cells.Item(8, 4).Value = 5.960285277117704
$ws.Cells.Item(8, 5).Value = 12.40768249843631
$ws.Cells.Item(8, 7).Value = 3.661654639550203
$ws.Cells.Item(8, 11).Value = 8.652743985268584
$ws.Cells.Item(8, 13).Value = 13.73572344235592
$ws.Cells.Item(8, 15).Value = 26.49013336761151

$ws.Cells.Item(9, 2).Value = 9.872238876781459
$ws.Cells.Item(9, 3).Value = 5.136787487997692
$ws.Cells.Item(9, 4).Value = 6.251223236339361
$ws.Cells.Item(9, 5).Value = 12.97889963735552
$ws.Cells.Item(9, 7).Value = 3.656362089714061
$ws.Cells.Item(9, 11).Value = 9.022111187644203
$ws.Cells.Item(9, 13).Value = 14.12367093014437
$ws.Cells.Item(9, 15).Value = 26.50505196234593

$ws.Cells.Item(10, 2).Value = 10.28566956315359
$ws.Cells.Item(10, 3).Value = 5.347079679399152
$ws.Cells.Item(10, 4).Value = 6.464393887474391
$ws.Cells.Item(10, 5).Value = 13.40174588825505
$ws.Cells.Item(10, 7).Value = 3.652823809355272
$ws.Cells.Item(10, 11).Value = 9.293191768040012
$ws.Cells.Item(10, 13).Value = 14.41886444244845
$ws.Cells.Item(10, 15).Value = 26.54287718086492

$ws.Cells.Item(11, 2).Value = 10.47045472539623
$ws.Cells.Item(11, 3).Value = 5.439101627760626
$ws.Cells.Item(11, 4).Value = 6.56074735578309
$ws.Cells.Item(11, 5).Value = 13.59384168954667
$ws.Cells.Item(11, 7).Value = 3.651289345374348
$ws.Cells.Item(11, 11).Value = 9.415856449484469
$ws.Cells.Item(11, 13).Value = 14.55480702675535
$ws.Cells.Item(11, 15).Value = 26.56590892279415

$ws.Cells.Item(12, 2).Value = 10.53988254407845
$ws.Cells.Item(12, 3).Value = 5.473408364019992
$ws.Cells.Item(12, 4).Value = 6.597105929160434
$ws.Cells.Item(12, 5).Value = 13.66646992066997
$ws.Cells.Item(12, 7).Value = 3.650719021613646
$ws.Cells.Item(12, 11).Value = 9.462165661033461
$ws.Cells.Item(12, 13).Value = 14.60647653109102
$ws.Cells.Item(12, 15).Value = 26.57546590821316

$ws.Cells.Item(13, 2).Value = 10.52495542888093
$ws.Cells.Item(13, 3).Value = 5.466044053229101
$ws.Cells.Item(13, 4).Value = 6.589281775229226
$ws.Cells.Item(13, 5).Value = 13.65083440509421
$ws.Cells.Item(13, 7).Value = 3.650841374086523
$ws.Cells.Item(13, 11).Value = 9.452199160673729
$ws.Cells.Item(13, 13).Value = 14.59534083442504
$ws.Cells.Item(13, 15).Value = 26.57337053104806

$ws.Cells.Item(14, 2).Value = 10.47617794378551
$ws.Cells.Item(14, 3).Value = 5.441934973438105
$ws.Cells.Item(14, 4).Value = 6.563741388932799
$ws.Cells.Item(14, 5).Value = 13.59981957751177
$ws.Cells.Item(14, 7).Value = 3.651242209487841
$ws.Cells.Item(14, 11).Value = 9.4196694319452
$ws.Cells.Item(14, 13).Value = 14.5590543501196
$ws.Cells.Item(14, 15).Value = 26.56667845531096

$ws.Cells.Item(15, 2).Value = 10.44622709375109
$ws.Cells.Item(15, 3).Value = 5.427096687284992
$ws.Cells.Item(15, 4).Value = 6.548079310520652
$ws.Cells.Item(15, 5).Value = 13.56855441060332
$ws.Cells.Item(15, 7).Value = 3.651489130292938
$ws.Cells.Item(15, 11).Value = 9.399724263599083
$ws.Cells.Item(15, 13).Value = 14.53685130447662
$ws.Cells.Item(15, 15).Value = 26.56268807540092

$ws.Cells.Item(16, 2).Value = 10.2735208476499
$ws.Cells.Item(16, 3).Value = 5.340991120196467
$ws.Cells.Item(16, 4).Value = 6.458081092932395
$ws.Cells.Item(16, 5).Value = 13.38918010382324
$ws.Cells.Item(16, 7).Value = 3.652925596829135
$ws.Cells.Item(16, 11).Value = 9.285158121981434
$ws.Cells.Item(16, 13).Value = 14.41000957440498
$ws.Cells.Item(16, 15).Value = 26.54148903886396

$ws.Cells.Item(17, 2).Value = 10.16667506241268
$ws.Cells.Item(17, 3).Value = 5.287223365896261
$ws.Cells.Item(17, 4).Value = 6.40268290379985
$ws.Cells.Item(17, 5).Value = 13.27901759712572
$ws.Cells.Item(17, 7).Value = 3.653826020978929
$ws.Cells.Item(17, 11).Value = 9.214674189787351
$ws.Cells.Item(17, 13).Value = 14.33258704810449
$ws.Cells.Item(17, 15).Value = 26.52997464591684

$ws.Cells.Item(18, 2).Value = 10.10491487804701
$ws.Cells.Item(18, 3).Value = 5.255955732241396
$ws.Cells.Item(18, 4).Value = 6.370763305392018
$ws.Cells.Item(18, 5).Value = 13.21563506676327
$ws.Cells.Item(18, 7).Value = 3.654350994951276
$ws.Cells.Item(18, 11).Value = 9.174075120693763
$ws.Cells.Item(18, 13).Value = 14.28821450328886
$ws.Cells.Item(18, 15).Value = 26.52390031685038

$ws.Cells.Item(19, 2).Value = 10.08395384637279
$ws.Cells.Item(19, 3).Value = 5.24531086504357
$ws.Cells.Item(19, 4).Value = 6.359947489254138
$ws.Cells.Item(19, 5).Value = 13.19417380108939
$ws.Cells.Item(19, 7).Value = 3.654529958807811
$ws.Cells.Item(19, 11).Value = 9.160320397282382
$ws.Cells.Item(19, 13).Value = 14.27321945127843
$ws.Cells.Item(19, 15).Value = 26.52193790548803

$ws.Cells.Item(20, 2).Value = 10.17808113043432
$ws.Cells.Item(20, 3).Value = 5.292982538952701
$ws.Cells.Item(20, 4).Value = 6.408586225923319
$ws.Cells.Item(20, 5).Value = 13.29074721924843
$ws.Cells.Item(20, 7).Value = 3.653729437593887
$ws.Cells.Item(20, 11).Value = 9.22218373845193
$ws.Cells.Item(20, 13).Value = 14.34081271241178
$ws.Cells.Item(20, 15).Value = 26.53114362448792

$ws.Cells.Item(21, 2).Value = 10.49052045822866
$ws.Cells.Item(21, 3).Value = 5.449031171676463
$ws.Cells.Item(21, 4).Value = 6.571247007682739
$ws.Cells.Item(21, 5).Value = 13.61480756784086
$ws.Cells.Item(21, 7).Value = 3.651124183267549
$ws.Cells.Item(21, 11).Value = 9.429228403832031
$ws.Cells.Item(21, 13).Value = 14.56970776222906
$ws.Cells.Item(21, 15).Value = 26.56862143022121

$ws.Cells.Item(22, 2).Value = 10.69150329127278
$ws.Cells.Item(22, 3).Value = 5.547864144554741
$ws.Cells.Item(22, 4).Value = 6.676790078306687
$ws.Cells.Item(22, 5).Value = 13.8259003277921
$ws.Cells.Item(22, 7).Value = 3.649484099483575
$ws.Cells.Item(22, 11).Value = 9.563700257338974
$ws.Cells.Item(22, 13).Value = 14.72039553091504
$ws.Cells.Item(22, 15).Value = 26.59798291564346

$ws.Cells.Item(23, 2).Value = 10.58455189945468
$ws.Cells.Item(23, 3).Value = 5.49540868369337
$ws.Cells.Item(23, 4).Value = 6.620542236616708
$ws.Cells.Item(23, 5).Value = 13.71332471406426
$ws.Cells.Item(23, 7).Value = 3.650353733963046
$ws.Cells.Item(23, 11).Value = 9.49202247378696
$ws.Cells.Item(23, 13).Value = 14.63988629170716
$ws.Cells.Item(23, 15).Value = 26.58186769052441

$ws.Cells.Item(24, 2).Value = 10.17292547900014
$ws.Cells.Item(24, 3).Value = 5.290379923942443
$ws.Cells.Item(24, 4).Value = 6.405917550405188
$ws.Cells.Item(24, 5).Value = 13.28544440271447
$ws.Cells.Item(24, 7).Value = 3.653773080172503
$ws.Cells.Item(24, 11).Value = 9.218788907423718
$ws.Cells.Item(24, 13).Value = 14.33709345393796
$ws.Cells.Item(24, 15).Value = 26.53061342998895

$ws.Cells.Item(25, 2).Value = 9.717597512225618
$ws.Cells.Item(25, 3).Value = 5.05623570326841
$ws.Cells.Item(25, 4).Value = 6.172444548156533
$ws.Cells.Item(25, 5).Value = 12.823480281092
$ws.Cells.Item(25, 7).Value = 3.657732089592313
$ws.Cells.Item(25, 11).Value = 8.922032286122828
$ws.Cells.Item(25, 13).Value = 14.01674699680888
$ws.Cells.Item(25, 15).Value = 26.49629730141043
